$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.351.08"
$ws.Range("E2").Value = "  -1.45%  "

$ws.Range("D3").Value = "'1.592.91"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  -0.54%  "

$ws.Range("D5").Value = "'209.98"

$ws.Range("D6").Value = "'0.507"
$ws.Range("E6").Value = "  -1.15%  "

$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("D8").Value = "'0.0612"
$ws.Range("E8").Value = "  -1.22%  "

$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("D10").Value = "'19.56"

$ws.Range("D11").Value = "'0.0844"
$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.08"
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.572.94"
$ws.Range("E14").Value = "  -1.96%  "

$ws.Range("E15").Value = "  -1.56%  "

$ws.Range("D16").Value = "'64.63"
$ws.Range("E16").Value = "  -0.77%  "

$ws.Range("D17").Value = "'26.345.79"
$ws.Range("E17").Value = "  -1.41%  "

$ws.Range("D18").Value = "'0.0₃0729"
$ws.Range("E18").Value = "  -2.03%  "

$ws.Range("D19").Value = "'7.48"
$ws.Range("E19").Value = "  +4.96%  "

$ws.Range("D20").Value = "'211.61"
$ws.Range("E20").Value = "  +0.86%  "

$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("E22").Value = "  -0.71%  "

$ws.Range("D23").Value = "'2.17"
$ws.Range("E23").Value = "  -3.32%  "

$ws.Range("D24").Value = "'8.93"
$ws.Range("E24").Value = "  -1.05%  "

$ws.Range("D25").Value = "'145.37"

$ws.Range("E26").Value = "  -0.42%  "

$ws.Range("E27").Value = "  -1.29%  "

$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").Value = "'15.30"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("E30").Value = "  -0.93%  "

$ws.Range("E31").Value = "  -0.65%  "

$ws.Range("E32").Value = "  -1.54%  "

$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").Value = "'1.303.50"
$ws.Range("E34").Value = "  +2.02%  "

$ws.Range("E35").Value = "  +2.93%  "

$ws.Range("E36").Value = "  -2.09%  "

$ws.Range("E37").Value = "  -0.87%  "

$ws.Range("E38").Value = "  -0.46%  "

$ws.Range("D39").Value = "'1.12"
$ws.Range("E39").Value = "  -12.14%  "

$ws.Range("E40").Value = "  -1.68%  "

$ws.Range("E42").Value = "  +3.17%  "

$ws.Range("D43").Value = "'62.74"
$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("E44").Value = "  -4.73%  "

$ws.Range("D45").Value = "'0.764"
$ws.Range("E45").Value = "  -1.90%  "

$ws.Range("D46").Value = "'1.728.74"
$ws.Range("E46").Value = "  -0.40%  "

$ws.Range("D47").Value = "'88.12"
$ws.Range("E47").Value = "  -2.65%  "

$ws.Range("E48").Value = "  -4.41%  "

$ws.Range("E49").Value = "  +8.38%  "

$ws.Range("D50").Value = "'0.0983"
$ws.Range("E50").Value = "  -4.55%  "

$ws.Range("E51").Value = "  -1.37%  "
